$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 70: section title "Resultados finais da análise", merged B70:C70 ---
$ws.Range("B70").Value = "Resultados finais da análise"
$ws.Range("B70:C70").Merge()

# --- Row 71-72: header row "Héuristica + buscas locais" / "Número de vezes que teve o melhor resultado" ---
$ws.Range("B71").Value = "Héuristica + buscas locais"
$ws.Range("C71").Value = "Número de vezes que teve o melhor resultado"
$ws.Range("B71:B72").Merge()
$ws.Range("C71:C72").Merge()
$ws.Rows("72").RowHeight = 26

# --- Rows 73-75: result counts ---
$ws.Range("B73").Value = "HVMP"
$ws.Range("C73").Formula = "=D65+J65+P65"

$ws.Range("B74").Value = "Híbrida 2"
$ws.Range("C74").Formula = "=F65+L65+R65"

$ws.Range("B75").Value = "Resultados iguais"
$ws.Range("C75").Formula = "=COUNTIF(F3:F64,""Igual"")+COUNTIF(L3:L64,""Igual"")+COUNTIF(R3:R64,""Igual"")"

# --- Borders + alignment for the new block (B70:C76) ---
$block = $ws.Range("B70:C75")
$block.Borders.LineStyle = 1
$block.Borders.Weight = 2
$block.HorizontalAlignment = -4108
$block.VerticalAlignment = -4108

$ws.Range("B71:C72").WrapText = $true

$ws.Cells.Item(76, 3).Value = $null
$ws.Cells.Item(76, 3).Borders.LineStyle = 1

# --- Column widths ---
$ws.Columns("B").ColumnWidth = 17.5
$ws.Columns("C").ColumnWidth = 20

# --- Match O65 / P65 / R65 formatting to N65 / Q65 (style reindex artifact in source edit) ---
$ws.Range("N65").Copy()
$ws.Range("O65:R65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View / selection state ---
$ws.Application.ActiveWindow.ScrollRow = 42
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B70:C75").Select()
